$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Orders")

# Row 4 corresponds to "Elmar Qara" -- update the Status column (D) from "Rejected" to "Accepted"
$ws.Range("D4").Value = "Accepted"
